# Remove the "demand2" row and the "net2" row from the element/type table.
# Before:
#   1 element | type
#   2 demand1 | demand
#   3 demand2 | demand   <- remove
#   4 net1    | net
#   5 net2    | net      <- remove
#   6 pv1     | pv
#   7 bat1    | bat
# After:
#   1 element | type
#   2 demand1 | demand
#   3 net1    | net
#   4 pv1     | pv
#   5 bat1    | bat

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (demand2 / demand).
$ws.Rows(3).Delete()

# After the above deletion, the former row 5 (net2 / net) is now row 4.
$ws.Rows(4).Delete()
